$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting rows 13-23 down to 14-24.
$ws.Rows(13).Insert()

# Copy the formatting/shared-string-typed boolean cells (G/H, "true") from the
# row above (row 12, which is also true/true) into the newly inserted row 13,
# so the cells keep the text-shared-string representation ("true") and the
# numFmt "@" style (style index 5) instead of being auto-typed as booleans.
$ws.Range("G12:H12").Copy()
$ws.Range("G13:H13").PasteSpecial()

# Fill in the rest of the new "角色" (Role) row data.
$ws.Range("A13").Value = 13
$ws.Range("B13").Value = "角色"
$ws.Range("C13").Value = "进入自己的角色"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 1
$ws.Range("N13").Value = "MainIcon10"

# Renumber the "序列" (Id) values for the two rows above the insertion point.
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# Grow the "表1" table to include the newly added row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N24"))

# Update the active selection shown when the sheet is next opened.
$ws.Range("B9").Select()
